# player fire moved into playerclass
# A new datasheet row (PLAYER_BULLET_RELOAD) is introduced right after
# PLAYER_BULLET_DELAY (old row 6), pushing every following row down by one.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new, blank row at row 6 - shifts rows 6..20 down to 7..21.
$ws.Rows("6").Insert()

# Populate the newly inserted row with the new player-fire-reload parameter.
$ws.Range("A6").Value = "PLAYER_BULLET_RELOAD"
$ws.Range("B6").Value = 10
$ws.Range("C6").Value = "int"
$ws.Range("D6").Value = "how many frames to wait to fire again"

# Match the author's final selection in the saved workbook.
$ws.Range("E9").Select()
